# Apply updated "想去人数" (want-to-go count) values across sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 0
$wsExpo.Range("F3").Value = 0
$wsExpo.Range("F6").Value = 515
$wsExpo.Range("F7").Value = 0
$wsExpo.Range("F8").Value = 232
$wsExpo.Range("F9").Value = 0
$wsExpo.Range("F10").Value = 77
$wsExpo.Range("F11").Value = 0

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 14
$wsShow.Range("F5").Value = 6

# --- Sheet "全部类型" (All Types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 0
$wsAll.Range("F3").Value = 0
$wsAll.Range("F6").Value = 0
$wsAll.Range("F7").Value = 0
$wsAll.Range("F8").Value = 515
$wsAll.Range("F10").Value = 0
$wsAll.Range("F11").Value = 232
$wsAll.Range("F12").Value = 511
$wsAll.Range("F13").Value = 77
$wsAll.Range("F15").Value = 3
$wsAll.Range("F16").Value = 0
$wsAll.Range("F17").Value = 0
